# "fix and upload course" — update the "course" row (row 7) of the task sheet:
#  - Start/Finish dates were recorded as free-text strings instead of real dates
#  - % Build finished (D7) bumped to 100%
#  - Assignee (E7) corrected from "Dev A" to "Nguyễn Trí Hậu"
#  - Leave the active selection on D7, matching where the edit was made

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "24/8/2022"
$ws.Range("C7").Value = "25/8/2022"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Nguyễn Trí Hậu"

$ws.Range("D7").Select()
